$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "Testmail #5: Wil je deze klant bellen?"
$ws.Range("B6").Value = "Beste collega,`nDank voor het doorsturen van deze testmail. Kun je wat meer context geven over wie deze klant is en waarom we hen zouden moeten bellen? Dan kan ik de juiste acties ondernemen.`nMet vriendelijke groet,`n[Jouw naam]  `n[Jouw functie]"
$ws.Range("C6").Value = "Wil je deze klant bellen?"
$ws.Range("D6").Value = "mailmind.test@zohomail.eu"
$ws.Range("E6").Value = "Klantenservice / Contact"
$ws.Range("F6").Value = "2025-08-04 20:18:33"
$ws.Range("G6").Value = "Ja"
$ws.Range("H6").Value = "Nee"
$ws.Range("I6").Value = "Ja"
$ws.Range("J6").Value = "Nee"
